$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Row 2: 240X120 PORCELANATO ---
$ws.Range("C2").Value2 = 129.6
$ws.Range("D2").Value2 = 0
$ws.Range("E2").Value2 = 129.6
$ws.Range("F2").Value2 = 0

# --- Row 3: 240X80 PORCELANATO ---
$ws.Range("C3").Value2 = 3592.51
$ws.Range("D3").Value2 = 0
$ws.Range("E3").Value2 = 3592.51
$ws.Range("F3").Value2 = 0

# --- Row 4: FREGADEROS DE COCINA ---
$ws.Range("C4").Value2 = 207.39
$ws.Range("D4").Value2 = 0
$ws.Range("E4").Value2 = 207.39
$ws.Range("F4").Value2 = 0

# --- Row 5: GRIFERIAS ---
$ws.Range("C5").Value2 = 86.41
$ws.Range("D5").Value2 = 0
$ws.Range("E5").Value2 = 86.41
$ws.Range("F5").Value2 = 0

# --- Row 6: INODOROS ---
$ws.Range("C6").Value2 = 660.6
$ws.Range("D6").Value2 = 23.4
$ws.Range("E6").Value2 = 637.2
$ws.Range("F6").Value2 = 0.03542234332425068

# --- Row 7: LAVABOS ---
$ws.Range("C7").Value2 = 93.90000000000001
$ws.Range("D7").Value2 = 0
$ws.Range("E7").Value2 = 93.90000000000001
$ws.Range("F7").Value2 = 0

# --- Row 8: NO RESURTIBLES ---
$ws.Range("C8").Value2 = 350
$ws.Range("D8").Value2 = 0
$ws.Range("E8").Value2 = 350
$ws.Range("F8").Value2 = 0

# --- Row 9: OTROS (unchanged) ---

# --- Row 10: PANELES DECORATIVOS ---
$ws.Range("D10").Value2 = 0
$ws.Range("E10").Value2 = 388.107983534392
$ws.Range("F10").Value2 = 0

# --- Row 11: PIEDRA SINTERIZADA ---
$ws.Range("C11").Value2 = 3446
$ws.Range("D11").Value2 = -142.56
$ws.Range("E11").Value2 = 3588.56
$ws.Range("F11").Value2 = -0.04136970400464306

# --- Row 12: PORCELANATO ---
$ws.Range("C12").Value2 = 31214
$ws.Range("D12").Value2 = 1592.89
$ws.Range("E12").Value2 = 29621.11
$ws.Range("F12").Value2 = 0.05103126802075992

# --- Row 13: PUERTAS DE SEGURIDAD (unchanged) ---

# --- Row 14: TOTAL ---
$ws.Range("C14").Value2 = 40279.56164865473
$ws.Range("D14").Value2 = 1473.73
$ws.Range("E14").Value2 = 38805.83164865473
$ws.Range("F14").Value2 = 0.03658753818760141

# --- Column width adjustments (D, E, F) ---
# Excel's ColumnWidth (character units) differs from the stored OOXML
# <col width="..."/> by a constant offset of 5/6 for this sheet's default
# font metrics, so subtract that offset to land exactly on the target
# stored width.
$ws.Columns.Item(4).ColumnWidth = 13 - (5/6)
$ws.Columns.Item(5).ColumnWidth = 22 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 26 - (5/6)
